$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (price + 1h volume change columns), per upstream scrape refresh.
# D-column "Price" values are stored as text (the source data uses dotted thousand
# separators like "27.585.13" that are not valid numbers) so we force each one with a
# leading apostrophe, matching how the original text-typed values behave in Excel, then
# reset the style back to Normal so no stray text-number-format is left on the cell.

$c = $ws.Range("D2")
$c.Value = '''27.585.13'
$c.Style = "Normal"
$ws.Range("E2").Value = '  +1.25%  '
$c = $ws.Range("D3")
$c.Value = '''1.760.25'
$c.Style = "Normal"
$ws.Range("E3").Value = '  -1.44%  '
$ws.Range("E4").Value = '  +0.01%  '
$c = $ws.Range("D5")
$c.Value = '''336.04'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.63%  '
$c = $ws.Range("D6")
$c.Value = '''1.001'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.03%  '
$c = $ws.Range("D7")
$c.Value = '''0.3829'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +0.81%  '
$c = $ws.Range("D8")
$c.Value = '''0.3398'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -0.54%  '
$c = $ws.Range("D9")
$c.Value = '''46.87'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -3.42%  '
$ws.Range("E10").Value = '  -5.06%  '
$c = $ws.Range("D11")
$c.Value = '''0.07368'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -1.61%  '
$c = $ws.Range("D12")
$c.Value = '''1.002'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -0.01%  '
$c = $ws.Range("D13")
$c.Value = '''22.28'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +1.96%  '
$c = $ws.Range("D14")
$c.Value = '''6.325'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -2.16%  '
$c = $ws.Range("D15")
$c.Value = '''1.759.61'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -1.57%  '
$c = $ws.Range("D16")
$c.Value = '''7.007'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -0.99%  '
$c = $ws.Range("D17")
$c.Value = '''0.00001073'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -2.08%  '
$c = $ws.Range("D18")
$c.Value = '''0.06650'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -0.05%  '
$c = $ws.Range("D19")
$c.Value = '''82.08'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -2.04%  '
$c = $ws.Range("D20")
$c.Value = '''1.002'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +0.04%  '
$c = $ws.Range("D21")
$c.Value = '''17.29'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -0.25%  '
$c = $ws.Range("D22")
$c.Value = '''6.372'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -3.84%  '
$c = $ws.Range("D23")
$c.Value = '''27.580.47'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +1.23%  '
$c = $ws.Range("D24")
$c.Value = '''12.03'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -2.74%  '
$c = $ws.Range("D25")
$c.Value = '''2.379'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -1.26%  '
$ws.Range("B26").Value = 'ImmutableX'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range("D26")
$c.Value = '''1.426'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -4.73%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range("D27")
$c.Value = '''20.60'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -3.40%  '
$c = $ws.Range("D28")
$c.Value = '''2.419'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -4.79%  '
$c = $ws.Range("D29")
$c.Value = '''152.43'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -0.36%  '
$c = $ws.Range("D30")
$c.Value = '''134.13'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -0.05%  '
$c = $ws.Range("D31")
$c.Value = '''1.959.34'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -1.51%  '
$c = $ws.Range("D32")
$c.Value = '''6.087'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +0.26%  '
$c = $ws.Range("D33")
$c.Value = '''3.963'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -1.42%  '
$c = $ws.Range("D34")
$c.Value = '''0.08793'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +1.08%  '
$c = $ws.Range("D35")
$c.Value = '''12.69'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -4.47%  '
$c = $ws.Range("D36")
$c.Value = '''0.02405'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +2.90%  '
$c = $ws.Range("D37")
$c.Value = '''0.6761'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -2.64%  '
$c = $ws.Range("D38")
$c.Value = '''5.315'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -2.42%  '
$c = $ws.Range("D39")
$c.Value = '''0.2175'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -1.43%  '
$c = $ws.Range("D40")
$c.Value = '''0.06267'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -1.02%  '
$c = $ws.Range("D41")
$c.Value = '''1.248'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +0.81%  '
$c = $ws.Range("D42")
$c.Value = '''1.505'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -9.03%  '
$c = $ws.Range("D43")
$c.Value = '''8.267'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -6.19%  '
$ws.Range("E44").Value = '  +0.04%  '
$c = $ws.Range("D45")
$c.Value = '''14.05'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -2.75%  '
$c = $ws.Range("D46")
$c.Value = '''0.6234'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -4.15%  '
$c = $ws.Range("D47")
$c.Value = '''3.824'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -0.43%  '
$c = $ws.Range("D48")
$c.Value = '''131.15'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +1.38%  '
$ws.Range("E49").Value = '  -3.82%  '
$c = $ws.Range("D50")
$c.Value = '''0.07346'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +3.07%  '
$c = $ws.Range("D51")
$c.Value = '''1.140'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +2.18%  '
